$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "Alexander BOLSHUNOV"
$ws.Cells.Item(2, 3).Value = 3482277
$ws.Cells.Item(2, 4).Value = 22366
$ws.Cells.Item(2, 5).Value = "m"
$ws.Cells.Item(2, 6).Value = 1569.883669528048
$ws.Cells.Item(2, 7).Value = 50
$ws.Cells.Item(2, 8).Value = 1

$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Andrey MELNICHENKO"
$ws.Cells.Item(3, 3).Value = 3481803
$ws.Cells.Item(3, 4).Value = 7488
$ws.Cells.Item(3, 5).Value = "m"
$ws.Cells.Item(3, 6).Value = 1535.813688753353
$ws.Cells.Item(3, 7).Value = 46
$ws.Cells.Item(3, 8).Value = 2

$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = "Dario COLOGNA"
$ws.Cells.Item(4, 3).Value = 3510023
$ws.Cells.Item(4, 4).Value = 4296
$ws.Cells.Item(4, 5).Value = "m"
$ws.Cells.Item(4, 6).Value = 1516.157717903772
$ws.Cells.Item(4, 7).Value = 43
$ws.Cells.Item(4, 8).Value = 3

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Denis SPITSOV"
$ws.Cells.Item(5, 3).Value = 3482280
$ws.Cells.Item(5, 4).Value = 2646
$ws.Cells.Item(5, 5).Value = "m"
$ws.Cells.Item(5, 6).Value = 1511.5749816709
$ws.Cells.Item(5, 7).Value = 40
$ws.Cells.Item(5, 8).Value = 4

$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Maurice MANIFICAT"
$ws.Cells.Item(6, 3).Value = 3190111
$ws.Cells.Item(6, 4).Value = 1450
$ws.Cells.Item(6, 5).Value = "m"
$ws.Cells.Item(6, 6).Value = 1511.276020310165
$ws.Cells.Item(6, 7).Value = 37
$ws.Cells.Item(6, 8).Value = 5

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Evgeniy BELOV"
$ws.Cells.Item(7, 3).Value = 3481161
$ws.Cells.Item(7, 4).Value = 3443
$ws.Cells.Item(7, 5).Value = "m"
$ws.Cells.Item(7, 6).Value = 1495.249811755924
$ws.Cells.Item(7, 7).Value = 34
$ws.Cells.Item(7, 8).Value = 6

$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 2).Value = "Clement PARISSE"
$ws.Cells.Item(8, 3).Value = 3190302
$ws.Cells.Item(8, 4).Value = 1470
$ws.Cells.Item(8, 5).Value = "m"
$ws.Cells.Item(8, 6).Value = 1492.45073227862
$ws.Cells.Item(8, 7).Value = 32
$ws.Cells.Item(8, 8).Value = 7

$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = "Ivan YAKIMUSHKIN"
$ws.Cells.Item(9, 3).Value = 3482105
$ws.Cells.Item(9, 4).Value = 5172
$ws.Cells.Item(9, 5).Value = "m"
$ws.Cells.Item(9, 6).Value = 1480.007133156578
$ws.Cells.Item(9, 7).Value = 30
$ws.Cells.Item(9, 8).Value = 8

$ws.Cells.Item(10, 1).Value = 33
$ws.Cells.Item(10, 2).Value = "Jean Marc GAILLARD"
$ws.Cells.Item(10, 3).Value = 1345875
$ws.Cells.Item(10, 4).Value = 900
$ws.Cells.Item(10, 5).Value = "m"
$ws.Cells.Item(10, 6).Value = 1466.446048751655
$ws.Cells.Item(10, 7).Value = 28
$ws.Cells.Item(10, 8).Value = 9

$ws.Cells.Item(11, 1).Value = 3
$ws.Cells.Item(11, 2).Value = "Artem MALTSEV"
$ws.Cells.Item(11, 3).Value = 3481432
$ws.Cells.Item(11, 4).Value = 4522
$ws.Cells.Item(11, 5).Value = "m"
$ws.Cells.Item(11, 6).Value = 1461.21214547662
$ws.Cells.Item(11, 7).Value = 26
$ws.Cells.Item(11, 8).Value = 10

$ws.Cells.Item(12, 1).Value = 18
$ws.Cells.Item(12, 2).Value = "Florian NOTZ"
$ws.Cells.Item(12, 3).Value = 3200376
$ws.Cells.Item(12, 4).Value = 753
$ws.Cells.Item(12, 5).Value = "m"
$ws.Cells.Item(12, 6).Value = 1455.205922000278
$ws.Cells.Item(12, 7).Value = 24
$ws.Cells.Item(12, 8).Value = 11

$ws.Cells.Item(13, 1).Value = 14
$ws.Cells.Item(13, 2).Value = "Adrien BACKSCHEIDER"
$ws.Cells.Item(13, 3).Value = 3190268
$ws.Cells.Item(13, 4).Value = 861
$ws.Cells.Item(13, 5).Value = "m"
$ws.Cells.Item(13, 6).Value = 1446.412751980616
$ws.Cells.Item(13, 7).Value = 22
$ws.Cells.Item(13, 8).Value = 12

$ws.Cells.Item(14, 1).Value = 15
$ws.Cells.Item(14, 2).Value = "Lucas BOEGL"
$ws.Cells.Item(14, 3).Value = 3200205
$ws.Cells.Item(14, 4).Value = 1319
$ws.Cells.Item(14, 5).Value = "m"
$ws.Cells.Item(14, 6).Value = 1437.352752882432
$ws.Cells.Item(14, 7).Value = 20
$ws.Cells.Item(14, 8).Value = 13

$ws.Cells.Item(15, 1).Value = 9
$ws.Cells.Item(15, 2).Value = "Francesco DE FABIANI"
$ws.Cells.Item(15, 3).Value = 3290379
$ws.Cells.Item(15, 4).Value = 3407
$ws.Cells.Item(15, 5).Value = "m"
$ws.Cells.Item(15, 6).Value = 1436.865066444981
$ws.Cells.Item(15, 7).Value = 18
$ws.Cells.Item(15, 8).Value = 14

$ws.Cells.Item(16, 1).Value = 16
$ws.Cells.Item(16, 2).Value = "Alexey CHERVOTKIN"
$ws.Cells.Item(16, 3).Value = 3482119
$ws.Cells.Item(16, 4).Value = 4000
$ws.Cells.Item(16, 5).Value = "m"
$ws.Cells.Item(16, 6).Value = 1422.425606714776
$ws.Cells.Item(16, 7).Value = 16
$ws.Cells.Item(16, 8).Value = 15

$ws.Cells.Item(17, 1).Value = 20
$ws.Cells.Item(17, 2).Value = "Jonas DOBLER"
$ws.Cells.Item(17, 3).Value = 3200356
$ws.Cells.Item(17, 4).Value = 579
$ws.Cells.Item(17, 5).Value = "m"
$ws.Cells.Item(17, 6).Value = 1417.933863845853
$ws.Cells.Item(17, 7).Value = 15
$ws.Cells.Item(17, 8).Value = 16

$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Hugo LAPALUS"
$ws.Cells.Item(18, 3).Value = 3190529
$ws.Cells.Item(18, 4).Value = 471
$ws.Cells.Item(18, 5).Value = "m"
$ws.Cells.Item(18, 6).Value = 1411.581535067255
$ws.Cells.Item(18, 7).Value = 14
$ws.Cells.Item(18, 8).Value = 17

$ws.Cells.Item(19, 1).Value = 25
$ws.Cells.Item(19, 2).Value = "Ireneu ESTEVE ALTIMIRAS"
$ws.Cells.Item(19, 3).Value = 3020003
$ws.Cells.Item(19, 4).Value = 577
$ws.Cells.Item(19, 5).Value = "m"
$ws.Cells.Item(19, 6).Value = 1395.752793772332
$ws.Cells.Item(19, 7).Value = 13
$ws.Cells.Item(19, 8).Value = 18

$ws.Cells.Item(20, 1).Value = 22
$ws.Cells.Item(20, 2).Value = "Jules LAPIERRE"
$ws.Cells.Item(20, 3).Value = 3190398
$ws.Cells.Item(20, 4).Value = 585
$ws.Cells.Item(20, 5).Value = "m"
$ws.Cells.Item(20, 6).Value = 1387.055989107855
$ws.Cells.Item(20, 7).Value = 12
$ws.Cells.Item(20, 8).Value = 19

$ws.Cells.Item(21, 1).Value = 35
$ws.Cells.Item(21, 2).Value = "Giandomenico SALVADORI"
$ws.Cells.Item(21, 3).Value = 3290407
$ws.Cells.Item(21, 4).Value = 432
$ws.Cells.Item(21, 5).Value = "m"
$ws.Cells.Item(21, 6).Value = 1371.734128859096
$ws.Cells.Item(21, 7).Value = 11
$ws.Cells.Item(21, 8).Value = 20

$ws.Cells.Item(22, 1).Value = 24
$ws.Cells.Item(22, 2).Value = "Beda KLEE"
$ws.Cells.Item(22, 3).Value = 3510534
$ws.Cells.Item(22, 4).Value = 309
$ws.Cells.Item(22, 5).Value = "m"
$ws.Cells.Item(22, 6).Value = 1365.995463814761
$ws.Cells.Item(22, 7).Value = 10
$ws.Cells.Item(22, 8).Value = 21

$ws.Cells.Item(23, 1).Value = 13
$ws.Cells.Item(23, 2).Value = "Federico PELLEGRINO"
$ws.Cells.Item(23, 3).Value = 3290326
$ws.Cells.Item(23, 4).Value = 11514
$ws.Cells.Item(23, 5).Value = "m"
$ws.Cells.Item(23, 6).Value = 1359.051059222495
$ws.Cells.Item(23, 7).Value = 9
$ws.Cells.Item(23, 8).Value = 22

$ws.Cells.Item(24, 1).Value = 27
$ws.Cells.Item(24, 2).Value = "Jonas BAUMANN"
$ws.Cells.Item(24, 3).Value = 3510342
$ws.Cells.Item(24, 4).Value = 627
$ws.Cells.Item(24, 5).Value = "m"
$ws.Cells.Item(24, 6).Value = 1351.949637974542
$ws.Cells.Item(24, 7).Value = 8
$ws.Cells.Item(24, 8).Value = 23

$ws.Cells.Item(25, 1).Value = 34
$ws.Cells.Item(25, 2).Value = "Thomas BING"
$ws.Cells.Item(25, 3).Value = 3200241
$ws.Cells.Item(25, 4).Value = 215
$ws.Cells.Item(25, 5).Value = "m"
$ws.Cells.Item(25, 6).Value = 1351.349482092891
$ws.Cells.Item(25, 7).Value = 7
$ws.Cells.Item(25, 8).Value = 24

$ws.Cells.Item(26, 1).Value = 30
$ws.Cells.Item(26, 2).Value = "Roman FURGER"
$ws.Cells.Item(26, 3).Value = 3510351
$ws.Cells.Item(26, 4).Value = 1000
$ws.Cells.Item(26, 5).Value = "m"
$ws.Cells.Item(26, 6).Value = 1348.727760884097
$ws.Cells.Item(26, 7).Value = 6
$ws.Cells.Item(26, 8).Value = 25

$ws.Cells.Item(27, 1).Value = 23
$ws.Cells.Item(27, 2).Value = "Michal NOVAK"
$ws.Cells.Item(27, 3).Value = 3150570
$ws.Cells.Item(27, 4).Value = 316
$ws.Cells.Item(27, 5).Value = "m"
$ws.Cells.Item(27, 6).Value = 1347.717348029758
$ws.Cells.Item(27, 7).Value = 5
$ws.Cells.Item(27, 8).Value = 26

$ws.Cells.Item(28, 1).Value = 17
$ws.Cells.Item(28, 2).Value = "Gus SCHUMACHER"
$ws.Cells.Item(28, 3).Value = 3530882
$ws.Cells.Item(28, 4).Value = 592
$ws.Cells.Item(28, 5).Value = "m"
$ws.Cells.Item(28, 6).Value = 1342.495273545723
$ws.Cells.Item(28, 7).Value = 4
$ws.Cells.Item(28, 8).Value = 27

$ws.Cells.Item(29, 1).Value = 21
$ws.Cells.Item(29, 2).Value = "William POROMAA"
$ws.Cells.Item(29, 3).Value = 3501741
$ws.Cells.Item(29, 4).Value = 990
$ws.Cells.Item(29, 5).Value = "m"
$ws.Cells.Item(29, 6).Value = 1323.597381206251
$ws.Cells.Item(29, 7).Value = 3
$ws.Cells.Item(29, 8).Value = 28

$ws.Cells.Item(30, 1).Value = 37
$ws.Cells.Item(30, 2).Value = "Scott PATTERSON"
$ws.Cells.Item(30, 3).Value = 3530532
$ws.Cells.Item(30, 4).Value = 258
$ws.Cells.Item(30, 5).Value = "m"
$ws.Cells.Item(30, 6).Value = 1323.225621873354
$ws.Cells.Item(30, 7).Value = 2
$ws.Cells.Item(30, 8).Value = 29

$ws.Cells.Item(31, 1).Value = 28
$ws.Cells.Item(31, 2).Value = "Janosch BRUGGER"
$ws.Cells.Item(31, 3).Value = 3200676
$ws.Cells.Item(31, 4).Value = 225
$ws.Cells.Item(31, 5).Value = "m"
$ws.Cells.Item(31, 6).Value = 1322.407345557064
$ws.Cells.Item(31, 7).Value = 1
$ws.Cells.Item(31, 8).Value = 30

$ws.Cells.Item(32, 1).Value = 52
$ws.Cells.Item(32, 2).Value = "Jessie DIGGINS"
$ws.Cells.Item(32, 3).Value = 3535410
$ws.Cells.Item(32, 4).Value = 12712
$ws.Cells.Item(32, 5).Value = "f"
$ws.Cells.Item(32, 6).Value = 1540.719084124432
$ws.Cells.Item(32, 7).Value = 50
$ws.Cells.Item(32, 8).Value = 1

$ws.Cells.Item(33, 1).Value = 54
$ws.Cells.Item(33, 2).Value = "Krista PARMAKOSKI"
$ws.Cells.Item(33, 3).Value = 3185256
$ws.Cells.Item(33, 4).Value = 5201
$ws.Cells.Item(33, 5).Value = "f"
$ws.Cells.Item(33, 6).Value = 1517.541323329707
$ws.Cells.Item(33, 7).Value = 46
$ws.Cells.Item(33, 8).Value = 2

$ws.Cells.Item(34, 1).Value = 57
$ws.Cells.Item(34, 2).Value = "Ebba ANDERSSON"
$ws.Cells.Item(34, 3).Value = 3505990
$ws.Cells.Item(34, 4).Value = 15649
$ws.Cells.Item(34, 5).Value = "f"
$ws.Cells.Item(34, 6).Value = 1501.218963683138
$ws.Cells.Item(34, 7).Value = 43
$ws.Cells.Item(34, 8).Value = 3

$ws.Cells.Item(35, 1).Value = 63
$ws.Cells.Item(35, 2).Value = "Teresa STADLOBER"
$ws.Cells.Item(35, 3).Value = 3055067
$ws.Cells.Item(35, 4).Value = 6614
$ws.Cells.Item(35, 5).Value = "f"
$ws.Cells.Item(35, 6).Value = 1492.768075175097
$ws.Cells.Item(35, 7).Value = 40
$ws.Cells.Item(35, 8).Value = 4

$ws.Cells.Item(36, 1).Value = 58
$ws.Cells.Item(36, 2).Value = "Rosie BRENNAN"
$ws.Cells.Item(36, 3).Value = 3535316
$ws.Cells.Item(36, 4).Value = 10007
$ws.Cells.Item(36, 5).Value = "f"
$ws.Cells.Item(36, 6).Value = 1471.853871681006
$ws.Cells.Item(36, 7).Value = 37
$ws.Cells.Item(36, 8).Value = 5

$ws.Cells.Item(37, 1).Value = 55
$ws.Cells.Item(37, 2).Value = "Natalia NEPRYAEVA"
$ws.Cells.Item(37, 3).Value = 3486010
$ws.Cells.Item(37, 4).Value = 16746
$ws.Cells.Item(37, 5).Value = "f"
$ws.Cells.Item(37, 6).Value = 1471.071379176677
$ws.Cells.Item(37, 7).Value = 34
$ws.Cells.Item(37, 8).Value = 6

$ws.Cells.Item(38, 1).Value = 53
$ws.Cells.Item(38, 2).Value = "Yulia STUPAK"
$ws.Cells.Item(38, 3).Value = 3486239
$ws.Cells.Item(38, 4).Value = 7809
$ws.Cells.Item(38, 5).Value = "f"
$ws.Cells.Item(38, 6).Value = 1440.030129025987
$ws.Cells.Item(38, 7).Value = 32
$ws.Cells.Item(38, 8).Value = 7

$ws.Cells.Item(39, 1).Value = 59
$ws.Cells.Item(39, 2).Value = "Katharina HENNIG"
$ws.Cells.Item(39, 3).Value = 3205460
$ws.Cells.Item(39, 4).Value = 4757
$ws.Cells.Item(39, 5).Value = "f"
$ws.Cells.Item(39, 6).Value = 1409.99822952098
$ws.Cells.Item(39, 7).Value = 30
$ws.Cells.Item(39, 8).Value = 8

$ws.Cells.Item(40, 1).Value = 67
$ws.Cells.Item(40, 2).Value = "Katerina RAZYMOVA"
$ws.Cells.Item(40, 3).Value = 3155314
$ws.Cells.Item(40, 4).Value = 1975
$ws.Cells.Item(40, 5).Value = "f"
$ws.Cells.Item(40, 6).Value = 1392.275202851005
$ws.Cells.Item(40, 7).Value = 28
$ws.Cells.Item(40, 8).Value = 9

$ws.Cells.Item(41, 1).Value = 70
$ws.Cells.Item(41, 2).Value = "Delphine CLAUDEL"
$ws.Cells.Item(41, 3).Value = 3195219
$ws.Cells.Item(41, 4).Value = 965
$ws.Cells.Item(41, 5).Value = "f"
$ws.Cells.Item(41, 6).Value = 1383.810372615237
$ws.Cells.Item(41, 7).Value = 26
$ws.Cells.Item(41, 8).Value = 10

$ws.Cells.Item(42, 1).Value = 66
$ws.Cells.Item(42, 2).Value = "Anamarija LAMPIC"
$ws.Cells.Item(42, 3).Value = 3565062
$ws.Cells.Item(42, 4).Value = 12607
$ws.Cells.Item(42, 5).Value = "f"
$ws.Cells.Item(42, 6).Value = 1383.36888670827
$ws.Cells.Item(42, 7).Value = 24
$ws.Cells.Item(42, 8).Value = 11

$ws.Cells.Item(43, 1).Value = 74
$ws.Cells.Item(43, 2).Value = "Jonna SUNDLING"
$ws.Cells.Item(43, 3).Value = 3505809
$ws.Cells.Item(43, 4).Value = 4518
$ws.Cells.Item(43, 5).Value = "f"
$ws.Cells.Item(43, 6).Value = 1378.953405844086
$ws.Cells.Item(43, 7).Value = 22
$ws.Cells.Item(43, 8).Value = 12

$ws.Cells.Item(44, 1).Value = 76
$ws.Cells.Item(44, 2).Value = "Anna NECHAEVSKAYA"
$ws.Cells.Item(44, 3).Value = 3485849
$ws.Cells.Item(44, 4).Value = 454
$ws.Cells.Item(44, 5).Value = "f"
$ws.Cells.Item(44, 6).Value = 1377.4252451641
$ws.Cells.Item(44, 7).Value = 20
$ws.Cells.Item(44, 8).Value = 13

$ws.Cells.Item(45, 1).Value = 56
$ws.Cells.Item(45, 2).Value = "Tatiana SORINA"
$ws.Cells.Item(45, 3).Value = 3485933
$ws.Cells.Item(45, 4).Value = 7690
$ws.Cells.Item(45, 5).Value = "f"
$ws.Cells.Item(45, 6).Value = 1366.876069063357
$ws.Cells.Item(45, 7).Value = 18
$ws.Cells.Item(45, 8).Value = 14

$ws.Cells.Item(46, 1).Value = 62
$ws.Cells.Item(46, 2).Value = "Nadine FAEHNDRICH"
$ws.Cells.Item(46, 3).Value = 3515221
$ws.Cells.Item(46, 4).Value = 10716
$ws.Cells.Item(46, 5).Value = "f"
$ws.Cells.Item(46, 6).Value = 1363.051778910941
$ws.Cells.Item(46, 7).Value = 16
$ws.Cells.Item(46, 8).Value = 15

$ws.Cells.Item(47, 1).Value = 72
$ws.Cells.Item(47, 2).Value = "Emma RIBOM"
$ws.Cells.Item(47, 3).Value = 3506008
$ws.Cells.Item(47, 4).Value = 2750
$ws.Cells.Item(47, 5).Value = "f"
$ws.Cells.Item(47, 6).Value = 1362.134781813197
$ws.Cells.Item(47, 7).Value = 15
$ws.Cells.Item(47, 8).Value = 16

$ws.Cells.Item(48, 1).Value = 61
$ws.Cells.Item(48, 2).Value = "Alisa ZHAMBALOVA"
$ws.Cells.Item(48, 3).Value = 3485862
$ws.Cells.Item(48, 4).Value = 1015
$ws.Cells.Item(48, 5).Value = "f"
$ws.Cells.Item(48, 6).Value = 1358.773969960442
$ws.Cells.Item(48, 7).Value = 14
$ws.Cells.Item(48, 8).Value = 17

$ws.Cells.Item(49, 1).Value = 71
$ws.Cells.Item(49, 2).Value = "Anna COMARELLA"
$ws.Cells.Item(49, 3).Value = 3295322
$ws.Cells.Item(49, 4).Value = 429
$ws.Cells.Item(49, 5).Value = "f"
$ws.Cells.Item(49, 6).Value = 1353.293399922565
$ws.Cells.Item(49, 7).Value = 13
$ws.Cells.Item(49, 8).Value = 18

$ws.Cells.Item(50, 1).Value = 69
$ws.Cells.Item(50, 2).Value = "Hailey SWIRBUL"
$ws.Cells.Item(50, 3).Value = 3535602
$ws.Cells.Item(50, 4).Value = 1197
$ws.Cells.Item(50, 5).Value = "f"
$ws.Cells.Item(50, 6).Value = 1352.073118082945
$ws.Cells.Item(50, 7).Value = 12
$ws.Cells.Item(50, 8).Value = 19

$ws.Cells.Item(51, 1).Value = 82
$ws.Cells.Item(51, 2).Value = "Moa LUNDGREN"
$ws.Cells.Item(51, 3).Value = 3506079
$ws.Cells.Item(51, 4).Value = 1793
$ws.Cells.Item(51, 5).Value = "f"
$ws.Cells.Item(51, 6).Value = 1344.767748996303
$ws.Cells.Item(51, 7).Value = 11
$ws.Cells.Item(51, 8).Value = 20

$ws.Cells.Item(52, 1).Value = 65
$ws.Cells.Item(52, 2).Value = "Maja DAHLQVIST"
$ws.Cells.Item(52, 3).Value = 3505800
$ws.Cells.Item(52, 4).Value = 4816
$ws.Cells.Item(52, 5).Value = "f"
$ws.Cells.Item(52, 6).Value = 1338.091514371696
$ws.Cells.Item(52, 7).Value = 10
$ws.Cells.Item(52, 8).Value = 21

$ws.Cells.Item(53, 1).Value = 64
$ws.Cells.Item(53, 2).Value = "Yana KIRPICHENKO"
$ws.Cells.Item(53, 3).Value = 3486314
$ws.Cells.Item(53, 4).Value = 730
$ws.Cells.Item(53, 5).Value = "f"
$ws.Cells.Item(53, 6).Value = 1337.699673860507
$ws.Cells.Item(53, 7).Value = 9
$ws.Cells.Item(53, 8).Value = 22

$ws.Cells.Item(54, 1).Value = 78
$ws.Cells.Item(54, 2).Value = "Pia FINK"
$ws.Cells.Item(54, 3).Value = 3205407
$ws.Cells.Item(54, 4).Value = 328
$ws.Cells.Item(54, 5).Value = "f"
$ws.Cells.Item(54, 6).Value = 1336.233588521824
$ws.Cells.Item(54, 7).Value = 8
$ws.Cells.Item(54, 8).Value = 23

$ws.Cells.Item(55, 1).Value = 60
$ws.Cells.Item(55, 2).Value = "Linn SVAHN"
$ws.Cells.Item(55, 3).Value = 3506166
$ws.Cells.Item(55, 4).Value = 4946
$ws.Cells.Item(55, 5).Value = "f"
$ws.Cells.Item(55, 6).Value = 1325.408854554127
$ws.Cells.Item(55, 7).Value = 7
$ws.Cells.Item(55, 8).Value = 24

$ws.Cells.Item(56, 1).Value = 77
$ws.Cells.Item(56, 2).Value = "Patricija EIDUKA"
$ws.Cells.Item(56, 3).Value = 3555052
$ws.Cells.Item(56, 4).Value = 368
$ws.Cells.Item(56, 5).Value = "f"
$ws.Cells.Item(56, 6).Value = 1315.507997812811
$ws.Cells.Item(56, 7).Value = 6
$ws.Cells.Item(56, 8).Value = 25

$ws.Cells.Item(57, 1).Value = 75
$ws.Cells.Item(57, 2).Value = "Hristina MATSOKINA"
$ws.Cells.Item(57, 3).Value = 3486563
$ws.Cells.Item(57, 4).Value = 540
$ws.Cells.Item(57, 5).Value = "f"
$ws.Cells.Item(57, 6).Value = 1300.731259692502
$ws.Cells.Item(57, 7).Value = 5
$ws.Cells.Item(57, 8).Value = 26

$ws.Cells.Item(58, 1).Value = 80
$ws.Cells.Item(58, 2).Value = "Izabela MARCISZ"
$ws.Cells.Item(58, 3).Value = 3435197
$ws.Cells.Item(58, 4).Value = 229
$ws.Cells.Item(58, 5).Value = "f"
$ws.Cells.Item(58, 6).Value = 1299.048453045386
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 27

$ws.Cells.Item(59, 1).Value = 86
$ws.Cells.Item(59, 2).Value = "Moa OLSSON"
$ws.Cells.Item(59, 3).Value = 3505998
$ws.Cells.Item(59, 4).Value = 227
$ws.Cells.Item(59, 5).Value = "f"
$ws.Cells.Item(59, 6).Value = 1296.041486984115
$ws.Cells.Item(59, 7).Value = 3
$ws.Cells.Item(59, 8).Value = 28

$ws.Cells.Item(60, 1).Value = 85
$ws.Cells.Item(60, 2).Value = "Caitlin PATTERSON"
$ws.Cells.Item(60, 3).Value = 3535385
$ws.Cells.Item(60, 4).Value = 211
$ws.Cells.Item(60, 5).Value = "f"
$ws.Cells.Item(60, 6).Value = 1295.733254937758
$ws.Cells.Item(60, 7).Value = 2
$ws.Cells.Item(60, 8).Value = 29

$ws.Cells.Item(61, 1).Value = 73
$ws.Cells.Item(61, 2).Value = "Katharine OGDEN"
$ws.Cells.Item(61, 3).Value = 3535601
$ws.Cells.Item(61, 4).Value = 217
$ws.Cells.Item(61, 5).Value = "f"
$ws.Cells.Item(61, 6).Value = 1289.357740969972
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 30
